$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status-check timestamp in F1
$ws.Range("F1").Value = "Last status check on: 09.02.2022 02:45"

# Row 9 (OMV IKEA): shift old price into C9, set new price in B9,
# write delta/date as plain text strings (use Value2 to dodge the
# broken Value getter in this runtime)
$ws.Range("C9").Value2 = $ws.Range("B9").Value2
$ws.Range("B9").Value2 = 38.5

# D9/E9 become plain text. A literal-string formula keeps the result
# as text without Excel's auto number/date parsing of "+0.6" / a date
# string, and (unlike NumberFormat="@") doesn't mint a brand new cell
# style. Copy/PasteSpecial values-only then collapses the formula back
# down to a plain literal, matching the target's non-formula text cell.
$ws.Range("D9").Formula = "=""+0.6"""
$ws.Range("D9").Copy($ws.Range("D9")) | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null

# E9 also loses its date number-format in the target, so clear the
# cell formatting back to the default style first.
$ws.Range("E9").ClearFormats()
$ws.Range("E9").Formula = "=""2022-02-09 02:46:44"""
$ws.Range("E9").Copy($ws.Range("E9")) | Out-Null
$ws.Range("E9").PasteSpecial(-4163) | Out-Null
